$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Helper: namespace-qualified "pkg:package" wrapper expected by
# Range.InsertXML (mirrors what real Word COM automation emits when it
# round-trips WordOpenXML through InsertXML).
# ---------------------------------------------------------------------
function New-PkgXml($innerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $innerBodyXml + '</w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) Title paragraph: "Week 7 Reading Guide Part 2: Hypothesis Tests"
#    gets re-typed as one run per word / inter-word space.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleStart = $titlePara.Range.Start
$titleEnd = $titlePara.Range.End
# Build a fresh Range via $d.Range(...) (not Paragraph.Range directly) -
# InsertXML reliably *replaces* such ranges, whereas a Paragraph.Range
# object (even after trimming its end) can make InsertXML *append*
# instead. Trim the trailing paragraph mark off the span so InsertXML
# only touches the visible text, preserving the paragraph's own mark.
$titleRange = $d.Range($titleStart, $titleEnd - 1)

$titleWords = @("Week", " ", "7", " ", "Reading", " ", "Guide", " ", "Part", " ", "2:", " ", "Hypothesis", " ", "Tests")
$titleRuns = ""
foreach ($w in $titleWords) {
    $titleRuns += '<w:r><w:t xml:space="preserve">' + $w + '</w:t></w:r>'
}
$titleXml = New-PkgXml('<w:p>' + $titleRuns + '</w:p>')
$titleRange.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2) The bold run containing the curly-quoted word "permute" is split
#    into three runs: opening quote / permute / closing quote, all
#    still bold, leaving the rest of the paragraph's runs untouched.
# ---------------------------------------------------------------------
$q1 = [char]0x201C
$q2 = [char]0x201D

$permutePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "*${q1}permute${q2}*") {
        $permutePara = $candidate
        break
    }
}

$permuteStart = $permutePara.Range.Start
$permuteEnd = $permutePara.Range.End
$permuteRange = $d.Range($permuteStart, $permuteEnd - 1)

$permuteRuns = '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">What does the word</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">' + $q1 + '</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">permute</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">' + $q2 + '</w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' + `
    '<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t xml:space="preserve">mean? How is this term related to the shuffling method that was used to obtain the slope statistic?</w:t></w:r>'

$permuteXml = New-PkgXml('<w:p><w:pPr><w:pStyle w:val="BodyText"/></w:pPr>' + $permuteRuns + '</w:p>')
$permuteRange.InsertXML($permuteXml)

# ---------------------------------------------------------------------
# 3) Style changes in styles.xml
# ---------------------------------------------------------------------

# Subtitle paragraph style now derives from Normal instead of Title.
$subtitleStyle = $d.Styles("Subtitle")
$subtitleStyle.BaseStyle = $d.Styles("Normal")

# Subtitle Char (the style linked to Subtitle) gains an explicit themed
# gray font color (text1, 65% tint -> 595959).
$subtitleCharStyle = $d.Styles("Subtitle")
$subtitleCharStyle.Font.ThemeColor = 1

# AbstractTitle paragraph style gains an explicit blue font color.
$abstractTitleStyle = $d.Styles("AbstractTitle")
$abstractTitleStyle.Font.Color = 9067060
